$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column G ("Recorded By") from 13 to 31 characters.
# Excel's ColumnWidth property adds ~0.8333 padding when round-tripped through
# the COM object model before it is stored back into the OOXML <col width="..."/>
# attribute, so we compensate by subtracting that offset to land exactly on 31.
$ws.Columns.Item(7).ColumnWidth = 30.1666667

# Populate the "Recorded By" column (G) for every session row that has already
# been recorded (these cells were previously blank placeholders).
$ws.Range("G2").Value = "Miss Dina Nasr"
$ws.Range("G4:G7").Value = "Miss Dina Nasr"
$ws.Range("G8:G12").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G13").Value = "Miss Dina Nasr"
$ws.Range("G14:G18").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G19:G20").Value = "Miss Dina Nasr"
$ws.Range("G21").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G22").Value = "Miss Dina Nasr"
$ws.Range("G23").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G24").Value = "Miss Dina Nasr"
$ws.Range("G28:G33").Value = "Miss Dina Nasr"
$ws.Range("G34:G38").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G39").Value = "Miss Dina Nasr"
$ws.Range("G40:G44").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G45:G46").Value = "Miss Dina Nasr"
$ws.Range("G47").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G48").Value = "Miss Dina Nasr"
$ws.Range("G49").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G50").Value = "Miss Dina Nasr"
$ws.Range("G54:G59").Value = "Miss Dina Nasr"
$ws.Range("G60:G64").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G65").Value = "Miss Dina Nasr"
$ws.Range("G66:G70").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G71:G72").Value = "Miss Dina Nasr"
$ws.Range("G73").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G74").Value = "Miss Dina Nasr"
$ws.Range("G75").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G76").Value = "Miss Dina Nasr"
$ws.Range("G80:G85").Value = "Miss Dina Nasr"
$ws.Range("G86:G90").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G91").Value = "Miss Dina Nasr"
$ws.Range("G92:G96").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G97:G98").Value = "Miss Dina Nasr"
$ws.Range("G99").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G100").Value = "Miss Dina Nasr"
$ws.Range("G101").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G102").Value = "Miss Dina Nasr"
$ws.Range("G106:G111").Value = "Miss Dina Nasr"
$ws.Range("G112:G116").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G117").Value = "Miss Dina Nasr"
$ws.Range("G118:G122").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G123:G124").Value = "Miss Dina Nasr"
$ws.Range("G125").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G126").Value = "Miss Dina Nasr"
$ws.Range("G127").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G128").Value = "Miss Dina Nasr"
$ws.Range("G133:G137").Value = "Miss Dina Nasr"
$ws.Range("G138:G142").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G143").Value = "Miss Dina Nasr"
$ws.Range("G144:G148").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G149:G150").Value = "Miss Dina Nasr"
$ws.Range("G151").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G152").Value = "Miss Dina Nasr"
$ws.Range("G153").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G154").Value = "Miss Dina Nasr"
$ws.Range("G158:G163").Value = "Miss Dina Nasr"
$ws.Range("G164").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G165:G166").Value = "Miss Dina Nasr"
$ws.Range("G167").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G168:G169").Value = "Miss Dina Nasr"
$ws.Range("G170").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G171:G173").Value = "Miss Dina Nasr"
$ws.Range("G174:G177").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G178:G180").Value = "Miss Dina Nasr"
$ws.Range("G185:G190").Value = "Miss Dina Nasr"
$ws.Range("G191").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G192:G193").Value = "Miss Dina Nasr"
$ws.Range("G194").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G195:G196").Value = "Miss Dina Nasr"
$ws.Range("G197").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G198:G200").Value = "Miss Dina Nasr"
$ws.Range("G201:G204").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G205:G207").Value = "Miss Dina Nasr"
$ws.Range("G212:G217").Value = "Miss Dina Nasr"
$ws.Range("G218").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G219:G220").Value = "Miss Dina Nasr"
$ws.Range("G221").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G222:G223").Value = "Miss Dina Nasr"
$ws.Range("G224").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G225:G227").Value = "Miss Dina Nasr"
$ws.Range("G228:G231").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G232:G234").Value = "Miss Dina Nasr"
$ws.Range("G239:G240").Value = "Miss Dina Nasr"
$ws.Range("G242:G244").Value = "Miss Dina Nasr"
$ws.Range("G245").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G246:G247").Value = "Miss Dina Nasr"
$ws.Range("G248").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G249:G250").Value = "Miss Dina Nasr"
$ws.Range("G251").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G252:G254").Value = "Miss Dina Nasr"
$ws.Range("G255:G258").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G259:G261").Value = "Miss Dina Nasr"
$ws.Range("G266:G271").Value = "Miss Dina Nasr"
$ws.Range("G272").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G273:G274").Value = "Miss Dina Nasr"
$ws.Range("G275").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G276:G277").Value = "Miss Dina Nasr"
$ws.Range("G278").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G279:G281").Value = "Miss Dina Nasr"
$ws.Range("G282:G285").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G286:G288").Value = "Miss Dina Nasr"
$ws.Range("G293:G298").Value = "Miss Dina Nasr"
$ws.Range("G299").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G300:G301").Value = "Miss Dina Nasr"
$ws.Range("G302").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G303:G304").Value = "Miss Dina Nasr"
$ws.Range("G305").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G306:G308").Value = "Miss Dina Nasr"
$ws.Range("G309:G312").Value = "Miss Dina Nasr, Administrator"
$ws.Range("G313:G315").Value = "Miss Dina Nasr"
